$d = $word.ActiveDocument

# 1. Update the date: "September 19, 2025" -> "September 21, 2025"
$r1 = $d.Content
$r1.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                  $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing-address paragraph:
#    "2930 Sanor Pl, Santa Clara CA 95051" ->
#       "2930 Sanor Pl" (existing paragraph)
#       "Santa Clara, CA 95051" (new paragraph, same run/paragraph formatting)
$r2 = $d.Content
$r2.Find.Execute(", Santa Clara CA 95051")
$r2.Text = "Santa Clara, CA 95051"
$r2.InsertParagraphBefore()

# 3. Remove the empty "NoSpacing" paragraph that follows "Board of Directors"
$r3 = $d.Content
$r3.Find.Execute("Board of Directors")
$para = $r3.Paragraphs(1)
$nextPara = $para.Next()
$nextPara.Range.Delete()
